# Update freelancer schedule cells to reflect reassigned shifts.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "C4"  = "15-24"
    "B10" = "7-16"
    "F10" = "15-24"
    "B11" = "7-16"
    "G12" = "7-16"
    "C13" = "10-19"
    "G13" = "off"
    "C14" = "10-19"
    "G14" = "7-16"
    "C16" = "15-24"
    "D16" = "7-16"
    "C17" = "10-19"
    "D17" = "15-24"
    "G17" = "off"
    "D18" = "7-16"
    "C19" = "10-19"
    "D19" = "15-24"
    "E19" = "15-24"
    "F19" = "off"
    "C20" = "10-19"
    "D20" = "15-24"
    "E20" = "15-24"
    "F20" = "off"
    "E23" = "7-16"
    "E24" = "7-16"
    "E25" = "7-16"
    "E26" = "7-16"
    "G26" = "15-24"
    "E27" = "7-16"
    "E28" = "7-16"
    "G29" = "7-16"
    "G30" = "7-16"
    "G31" = "7-16"
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
